# "Fix for square brackets in file names" -- adds a note about PowerShell 1.0
# test results (alongside the existing PowerShell 2.0 note) to the two cells
# that document the Get-Acl / GetAccessControl() workaround.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E17: keep the existing explanatory text, but insert a new line noting the
# workaround was also tested on PowerShell 1.0, right before the existing
# "PowerShell 2.0" line.
$ws.Range("E17").Value = "To get ACL, escape with grave accent, then do Get-Item. Then, when the item is returned, access its method GetAccessControl()`nTested on PowerShell 1.0 on Windows Server 2008 SP1`nTested on PowerShell 2.0 on Windows Server 2008 SP1"

# E18: same addition -- a "PowerShell 1.0" line ahead of the pre-existing
# "PowerShell 2.0" line.
$ws.Range("E18").Value = "Tested on PowerShell 1.0 on Windows Server 2008 SP1`nTested on PowerShell 2.0 on Windows Server 2008 SP1"

# The extra wrapped line makes row 17 taller (three lines -> four lines of
# wrapped text at s=4/wrapText), matching Excel's own auto re-measurement.
$ws.Rows.Item(17).RowHeight = 132
